$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (preserve formatting like trailing zeros / multi-dot numbers)
# by setting NumberFormat to Text before assigning, then resetting style so the
# cell format matches the original (unstyled) cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.771.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.907.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.908.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.96%  "
$ws.Range("E11").Value = "  +4.65%  "
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.388.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.770.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.909.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "435.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +20.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.699.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "341.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.90%  "
